$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "date coverage" labels in row 6 (shared across many columns) ---
$used = $ws.UsedRange
$numCols = $used.Columns.Count()
for ($c = 1; $c -le $numCols; $c++) {
    $cell = $ws.Cells.Item(6, $c)
    $v = $cell.Value()
    if ($v -eq "2006-03:2023-11") {
        $cell.Value = "2006-03:2023-12"
    } elseif ($v -eq "2006-07:2023-11") {
        $cell.Value = "2006-07:2023-12"
    } elseif ($v -eq "2018-02:2023-11") {
        $cell.Value = "2018-02:2023-12"
    }
}

# --- Update the "last updated" date in row 8 (same label repeated across columns) ---
# NOTE: these cells hold the date as literal text (quote-prefixed in the
# original file), so a leading apostrophe is required or Excel's COM layer
# will silently re-interpret the string as a real date serial number.
for ($c = 1; $c -le $numCols; $c++) {
    $cell = $ws.Cells.Item(8, $c)
    $v = $cell.Value()
    if ($v -eq "2023-12-29") {
        $cell.Value = "'2024-02-05"
    }
}

# --- Append a new data row (218) for 2023-12-31, copying number formats from row 217 ---
$ws.Range("A217:DV217").Copy()
$ws.Range("A218").PasteSpecial(-4122)
$excel.CutCopyMode = 0


$ws.Range("A218").Value = 45291
$ws.Range("B218").Value = 54695707.839000002
$ws.Range("C218").Value = 1382881.3624
$ws.Range("D218").Value = 441820.09529999999
$ws.Range("E218").Value = 15375.296399999999
$ws.Range("F218").Value = 413234.34570000001
$ws.Range("G218").Value = 175342.32870000001
$ws.Range("H218").Value = 337109.53700000001
$ws.Range("I218").Value = 194585.72159999999
$ws.Range("J218").Value = 2758049.9774000002
$ws.Range("K218").Value = 973687.44909999997
$ws.Range("L218").Value = 694536.11560000002
$ws.Range("M218").Value = 395663.76549999998
$ws.Range("N218").Value = 321984.83630000002
$ws.Range("O218").Value = 206551.6544
$ws.Range("P218").Value = 165626.10649999999
$ws.Range("Q218").Value = 40508005.241300002
$ws.Range("R218").Value = 751462.84889999998
$ws.Range("S218").Value = 498502.935
$ws.Range("T218").Value = 134849.87539999999
$ws.Range("U218").Value = 33784.384400000003
$ws.Range("V218").Value = 358901.0612
$ws.Range("W218").Value = 217924.62
$ws.Range("X218").Value = 379935.79200000002
$ws.Range("Y218").Value = 209004.8653
$ws.Range("Z218").Value = 1662010.253
$ws.Range("AA218").Value = 883838.0956
$ws.Range("AB218").Value = 143893.22880000001
$ws.Range("AC218").Value = 165839.70079999999
$ws.Range("AD218").Value = 22628.575499999999
$ws.Range("AE218").Value = 1707176.2416999999
$ws.Range("AF218").Value = 521107.29340000002
$ws.Range("AG218").Value = 5055372.4167999998
$ws.Range("AH218").Value = 615931.80850000004
$ws.Range("AI218").Value = 819551.90370000002
$ws.Range("AJ218").Value = 82968.9807
$ws.Range("AK218").Value = 591287.71160000004
$ws.Range("AL218").Value = 488248.09700000001
$ws.Range("AM218").Value = 41753.519500000002
$ws.Range("AN218").Value = 99324.510699999999
$ws.Range("AO218").Value = 437559.31199999998
$ws.Range("AP218").Value = 1578137.9339999999
$ws.Range("AQ218").Value = 383671.36989999999
$ws.Range("AR218").Value = 1194376.2287999999
$ws.Range("AS218").Value = 3562855.9662000001
$ws.Range("AT218").Value = 860840.04870000004
$ws.Range("AU218").Value = 212916.41469999999
$ws.Range("AV218").Value = 302559.86930000002
$ws.Range("AW218").Value = 88589.315600000002
$ws.Range("AX218").Value = 5642775.8417999996
$ws.Range("AY218").Value = 1517115.2651
$ws.Range("AZ218").Value = 3962050.6836000001
$ws.Range("BA218").Value = 6822935.0780999996
$ws.Range("BB218").Value = 5519045.5323000001
$ws.Range("BC218").Value = 200013.48740000001
$ws.Range("BD218").Value = 116971.22139999999
$ws.Range("BE218").Value = 2238872.1469999999
$ws.Range("BF218").Value = 642756.97219999996
$ws.Range("BG218").Value = 1104722.0355
$ws.Range("BH218").Value = 7114.0357999999997
$ws.Range("BI218").Value = 414018.83630000002
$ws.Range("BJ218").Value = 46667.370199999998
$ws.Range("BK218").Value = 747108.1361
$ws.Range("BL218").Value = 37330.992100000003
$ws.Range("BM218").Value = 319353.95689999999
$ws.Range("BN218").Value = 25691.700099999998
$ws.Range("BO218").Value = 9694.8269
$ws.Range("BP218").Value = 41043.368000000002
$ws.Range("BQ218").Value = 1720807.7662
$ws.Range("BR218").Value = 471967.18349999998
$ws.Range("BS218").Value = 2215731.8974000001
$ws.Range("BT218").Value = 38387.313199999997
$ws.Range("BU218").Value = 101991.1833
$ws.Range("BV218").Value = 76997.788
$ws.Range("BW218").Value = 680844.89080000005
$ws.Range("BX218").Value = 174591.98629999999
$ws.Range("BY218").Value = 79944.197400000005
$ws.Range("BZ218").Value = 11429653.566099999
$ws.Range("CA218").Value = 10523960.3693
$ws.Range("CB218").Value = 5284937.0719999997
$ws.Range("CC218").Value = 3041053.4630999998
$ws.Range("CD218").Value = 727859.29200000002
$ws.Range("CE218").Value = 192736.00450000001
$ws.Range("CF218").Value = 712956.18500000006
$ws.Range("CG218").Value = 1049594.2955
$ws.Range("CH218").Value = 360781.71659999999
$ws.Range("CI218").Value = 346086.50079999998
$ws.Range("CJ218").Value = 76851.416200000007
$ws.Range("CK218").Value = 265872.6421
$ws.Range("CL218").Value = 2092245.469
$ws.Range("CM218").Value = 1071386.4679
$ws.Range("CN218").Value = 414701.33539999998
$ws.Range("CO218").Value = 50919.480300000003
$ws.Range("CP218").Value = 3960.7977999999998
$ws.Range("CQ218").Value = 64679.915800000002
$ws.Range("CR218").Value = 28819.048999999999
$ws.Range("CS218").Value = 343394.92369999998
$ws.Range("CT218").Value = 83957.343500000003
$ws.Range("CU218").Value = 846657.85789999994
$ws.Range("CV218").Value = 34386.962399999997
$ws.Range("CW218").Value = 1272578.3188
$ws.Range("CX218").Value = 573588.41500000004
$ws.Range("CY218").Value = 499272.6471
$ws.Range("CZ218").Value = 251502.90359999999
$ws.Range("DA218").Value = 199717.26939999999
$ws.Range("DB218").Value = 3382335.9103000001
$ws.Range("DC218").Value = 604627.30709999998
$ws.Range("DD218").Value = 1029672.8101
$ws.Range("DE218").Value = 204696.25820000001
$ws.Range("DF218").Value = 1734103.3499
$ws.Range("DG218").Value = 773007.80350000004
$ws.Range("DH218").Value = 69221.586299999995
$ws.Range("DI218").Value = 4831967.4817000004
$ws.Range("DJ218").Value = 251382.11730000001
$ws.Range("DK218").Value = 40459.523699999998
$ws.Range("DL218").Value = 7501.3064999999997
$ws.Range("DM218").Value = 710156.27899999998
$ws.Range("DN218").Value = 82873.371400000004
$ws.Range("DO218").Value = 373101.52279999998
$ws.Range("DP218").Value = 1688548.5323000001
$ws.Range("DQ218").Value = 1460041.7254999999
$ws.Range("DR218").Value = 526125.4669
$ws.Range("DS218").Value = 710465.37040000001
$ws.Range("DT218").Value = 945288.97569999995
$ws.Range("DU218").Value = 13181065.8577
$ws.Range("DV218").Value = 85629856.122999996


Write-Output "done"
